# update statistics by excel
#
# The last data row (row 18: 2018-05-22 / 002611 / 博时黄金ETF联接C /
# 0.9326 / 0.9328) is copied and pasted into an 8-row x 2-column-block
# block below it (rows 21-28, columns A:E and again F:J), the same way a
# user would Copy a small range and Paste it repeatedly to fill a larger
# destination. This also overwrites the stray note that used to live in
# C22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A18:E18")
$source.Copy()

for ($r = 21; $r -le 28; $r++) {
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial()
    $ws.Range("F" + $r + ":J" + $r).PasteSpecial()
}

$excel.CutCopyMode = $false

# Leave the selection where the user ended up after the edit.
$ws.Range("O18").Select()
